# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-30 13:27:08
#
# This normalizes the "Recorded By" (column G) text on the
# "Session Analysis Results" sheet:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "system, System, backup@backdoor.com"
# Any other text in column G is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value()

    if ($null -eq $val) { continue }

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
